$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-28 (A=Active, B=Name, C=xDegrees, D=yDegrees, E=Population, F=AffectedPop, G=MaxDistance)
$data = @(
    @(2, 0, 'Balite', 14.8956, 120.7855, 5016, 0, 1000),
    @(3, 1, 'Balungao', 14.9143, 120.7622, 5720, 0, 1000),
    @(4, 1, 'Buguion', 14.894, 120.7985, 3841, 0, 1000),
    @(5, 1, 'Bulusan', 14.9076, 120.7455, 2603, 0, 1000),
    @(6, 1, 'Calizon', 14.9125, 120.753, 2221, 0, 1000),
    @(7, 1, 'Caniogan', 14.9054, 120.7733, 4510, 0, 1000),
    @(8, 1, 'Corazon', 14.9128, 120.7686, 2175, 0, 1000),
    @(9, 1, 'Frances', 14.9153, 120.7532, 6129, 6, 1000),
    @(10, 1, 'Gatbuca', 14.9218, 120.7685, 6384, 115, 1000),
    @(11, 1, 'Gugo', 14.9014, 120.7548, 1960, 57, 1000),
    @(12, 1, 'Iba Este', 14.8899, 120.7673, 4161, 0, 1000),
    @(13, 1, 'Iba O''Este', 14.8919, 120.7635, 14085, 601, 1000),
    @(14, 1, 'Meysulao', 14.9078, 120.7397, 4280, 56, 1000),
    @(15, 1, 'Meyto', 14.8831, 120.7295, 2925, 6, 1000),
    @(16, 1, 'Palimbang', 14.8994, 120.7756, 1684, 0, 1000),
    @(17, 1, 'Panducot', 14.8761, 120.738, 1752, 0, 1000),
    @(18, 1, 'Pio Cruzcosa', 14.8881, 120.7855, 4663, 92, 1000),
    @(19, 1, 'Poblacion', 14.9157, 120.7672, 1785, 0, 1000),
    @(20, 1, 'Pungo', 14.9023, 120.7914, 9528, 0, 1000),
    @(21, 1, 'San Jose', 14.8838, 120.7395, 5661, 0, 1000),
    @(22, 1, 'San Marcos', 14.8976, 120.7797, 2671, 0, 1000),
    @(23, 1, 'San Miguel', 14.917, 120.7427, 6005, 17, 1000),
    @(24, 1, 'Santa Lucia', 14.8982, 120.736, 2460, 0, 1000),
    @(25, 1, 'Santo Niño', 14.9047, 120.7792, 2544, 0, 1000),
    @(26, 1, 'Sapang Bayan', 14.9196, 120.7739, 3140, 65, 1000),
    @(27, 1, 'Sergio Bayan', 14.894, 120.7909, 1727, 0, 1000),
    @(28, 1, 'Sucol', 14.9138, 120.7701, 1059, 0, 1000)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = [bool]($row[1])
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
}

# Clear remark text that no longer applies
$ws.Cells.Item(6, 8).ClearContents()
$ws.Cells.Item(18, 8).ClearContents()
